# Update Mappings 22 Ontologies
# Adds a new "EDAM_DEF" column (F) with EDAM definitions for each mapped row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (copy the header formatting used by the other header cells)
$ws.Range("F1").Value2 = "EDAM_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# New EDAM_DEF values, one per data row (rows 2-10)
$defs = @(
    "['An array of numerical values.']",
    "['The concentration of a chemical compound.']",
    "['The spectrum of frequencies of electromagnetic radiation emitted from a molecule as a result of some spectroscopy experiment.']",
    "['The study of matter by studying the interference pattern from firing electrons at a sample, to analyse structures at resolutions higher than can be achieved using light.']",
    "['The study of matter and their structure by means of the diffraction of X-rays, typically the diffraction pattern caused by the regularly spaced atoms of a crystalline sample.']",
    "['An analytical chemistry technique that measures the mass-to-charge ratio and abundance of ions in the gas phase.']",
    "['Spectra from mass spectrometry.']",
    "['The study of matter by studying the diffraction pattern from firing neutrons at a sample, typically to determine atomic and/or magnetic structure.']",
    "['Construct some entity (typically a molecule sequence) from component pieces.']"
)

for ($i = 0; $i -lt $defs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value2 = $defs[$i]
}
